$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title ("Play Football: Champions Cup Free Slot
#    Game").
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2) Near the end of the document, right before the "Prompt for DALLE: ..."
#    paragraph, insert a new bold paragraph reading
#    "Play Football: Champions Cup Free Slot Game", matching the
#    <w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>...</w:t></w:r> run layout.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dallePara = $d.Paragraphs($count)
$insertionPoint = $d.Range($dallePara.Range.Start, $dallePara.Range.Start)

$xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Football: Champions Cup Free Slot Game</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xmlFrag)

# InsertXML leaves a spare empty paragraph between the new heading and
# the (former) DALLE paragraph -- drop it.
$count2 = $d.Paragraphs.Count
$spacerPara = $d.Paragraphs($count2 - 1)
$spacerPara.Range.Delete()

# ---------------------------------------------------------------------
# 3) Turn the old DALLE-prompt paragraph (still italic) into the new
#    meta-description blurb, keeping its italic run formatting. Using
#    $d.Range(start, end) (rather than the paragraph's own .Range)
#    replaces the text in place instead of inserting before it -- this
#    matters for the very last paragraph in the story.
# ---------------------------------------------------------------------
$count3 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($count3)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)
$finalRange.Text = "Read our review of Football: Champions Cup, a football-themed slot game with exciting bonus features and customizable betting options. Play for free now!"
